$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- Header row (row 1): new columns H=date, I=legislator_name, J=legislator_id ---
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

$header = $ws.Range("H1:J1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# --- Data rows (2-11): populate date / legislator_name / legislator_id for every stock row ---
for ($r = 2; $r -le 11; $r++) {
    $dateCell = $ws.Cells.Item($r, 8)
    $dateCell.Value = "'2012-02-01"
    $dateCell.ClearFormats()

    $ws.Cells.Item($r, 9).Value = "鄭汝芬"
    $ws.Cells.Item($r, 10).Value = 1713
}
